$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '27.589.16'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '
$ws.Cells.Item(3, 4).Value = '1.647.17'
$ws.Cells.Item(3, 5).Value = '  -0.47%  '
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '212.78'
$ws.Cells.Item(5, 5).Value = '  -0.49%  '
$ws.Cells.Item(6, 5).Value = '  +4.71%  '
$ws.Cells.Item(7, 5).Value = '  -0.09%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '23.56'
$ws.Cells.Item(8, 5).Value = '  -1.52%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.256'
$ws.Cells.Item(9, 5).Value = '  -1.66%  '
$ws.Cells.Item(10, 5).Value = '  -1.05%  '
$ws.Cells.Item(11, 5).Value = '  +1.28%  '
$ws.Cells.Item(12, 4).Value = '1.882.07'
$ws.Cells.Item(12, 5).Value = '  -0.46%  '
$ws.Cells.Item(13, 4).Value = '1.661.33'
$ws.Cells.Item(13, 5).Value = '  +0.38%  '
$ws.Cells.Item(14, 5).Value = '  +4.31%  '
$ws.Cells.Item(15, 5).Value = '  -2.13%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '64.43'
$ws.Cells.Item(16, 5).Value = '  -1.89%  '
$ws.Cells.Item(17, 4).Value = '27.553.15'
$ws.Cells.Item(17, 5).Value = '  +0.06%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '231.46'
$ws.Cells.Item(18, 5).Value = '  -3.39%  '
$ws.Cells.Item(19, 5).Value = '  -0.72%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.56'
$ws.Cells.Item(20, 5).Value = '  +0.38%  '
$ws.Cells.Item(21, 5).Value = '  -0.02%  '
$ws.Cells.Item(22, 5).Value = '  -2.80%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '9.75'
$ws.Cells.Item(23, 5).Value = '  +4.46%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.00'
$ws.Cells.Item(24, 5).Value = '  -1.75%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '149.25'
$ws.Cells.Item(25, 5).Value = '  +2.55%  '
$ws.Cells.Item(26, 5).Value = '  -2.57%  '
$ws.Cells.Item(27, 5).Value = '  +1.68%  '
$ws.Cells.Item(28, 5).Value = '  +0.08%  '
$ws.Cells.Item(29, 5).Value = '  -3.72%  '
$ws.Cells.Item(30, 5).Value = '  -2.04%  '
$ws.Cells.Item(31, 5).Value = '  -2.67%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.30'
$ws.Cells.Item(32, 5).Value = '  -0.23%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '3.19'
$ws.Cells.Item(33, 5).Value = '  +3.11%  '
$ws.Cells.Item(34, 4).Value = '1.430.62'
$ws.Cells.Item(34, 5).Value = '  -1.10%  '
$ws.Cells.Item(35, 5).Value = '  +2.97%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.37'
$ws.Cells.Item(36, 5).Value = '  -0.61%  '
$ws.Cells.Item(37, 5).Value = '  +0.30%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.883'
$ws.Cells.Item(38, 5).Value = '  -3.83%  '
$ws.Cells.Item(39, 5).Value = '  -2.30%  '
$ws.Cells.Item(40, 5).Value = '  -0.91%  '
$ws.Cells.Item(41, 5).Value = '  +0.01%  '
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.817'
$ws.Cells.Item(42, 5).Value = '  +2.86%  '
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '5.52'
$ws.Cells.Item(43, 5).Value = '  +1.98%  '
$ws.Cells.Item(44, 5).Value = '  +1.43%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '65.07'
$ws.Cells.Item(45, 5).Value = '  -5.59%  '
$ws.Cells.Item(46, 4).Value = '1.790.48'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.69'
$ws.Cells.Item(47, 5).Value = '  -0.79%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '88.06'
$ws.Cells.Item(48, 5).Value = '  -0.34%  '
$ws.Cells.Item(49, 5).Value = '  +0.75%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0998'
$ws.Cells.Item(50, 5).Value = '  -2.16%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '7.78'
$ws.Cells.Item(51, 5).Value = '  -0.56%  '
